# Adds a new "ODI Bowling Extra" worksheet (scraped extra bowling attributes)
# as the last sheet in the workbook, mirroring the layout of the existing
# "ODI Batting Extra" sheet but for bowling-related statistics.

$wb = $excel.ActiveWorkbook

# Create the new worksheet after the current last sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ODI Bowling Extra"

# Header row
$headers = @("MATCH_CODE", "MAIDEN_OVERS", "PERCENT_WICKETS_OF_ALL")
for ($col = 1; $col -le $headers.Length; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$col - 1]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
    $cell.Borders.LineStyle = 1
}

# Data rows: MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL
$rows = @(
    @("3712", "", ""),
    @("3714", "0", "20.00%"),
    @("3716", "0", "20.00%"),
    @("3718", "0", "20.00%"),
    @("3727", "0", "10.00%"),
    @("3780", "0", "20.00%"),
    @("3785", "2", "20.00%"),
    @("3804", "0", "10.00%"),
    @("3805", "0", "10.00%"),
    @("3885", "0", ""),
    @("3887", "0", ""),
    @("3889", "", ""),
    @("3908", "0", "10.00%"),
    @("3911", "0", ""),
    @("3930", "", ""),
    @("3932", "0", ""),
    @("4401", "", ""),
    @("4405", "", ""),
    @("4408", "", ""),
    @("4660", "", "")
)

$rowIndex = 2
foreach ($r in $rows) {
    $matchCode = $r[0]
    $maidenOvers = $r[1]
    $percentWickets = $r[2]

    $matchCodeCell = $ws.Cells.Item($rowIndex, 1)
    $matchCodeCell.NumberFormat = "@"
    $matchCodeCell.Value = $matchCode

    $maidenCell = $ws.Cells.Item($rowIndex, 2)
    $maidenCell.NumberFormat = "@"
    $maidenCell.Value = $maidenOvers

    $percentCell = $ws.Cells.Item($rowIndex, 3)
    $percentCell.NumberFormat = "@"
    $percentCell.Value = $percentWickets

    $rowIndex++
}

$ws.Range("A1").Select()
